$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "In Translation"
#    This shared string is used by the "Status" column on every sheet
#    (Overview!E2:F2, zh-cn!C2, de-de!C2). Note: comparing with the literal
#    on the LEFT of -eq avoids PowerShell coercing a boolean cell value
#    (e.g. the "To be localized" column, which holds True/False) into a
#    string-truthiness match.
# ---------------------------------------------------------------------------
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        if ("Ready for handoff" -eq $cell.Value2) {
            $cell.Value = "In Translation"
        }
    }
}

# ---------------------------------------------------------------------------
# 2. Narrow the Status-adjacent columns.
#    Overview columns E (zh-cn) & F (de-de), and column C ("Status") on the
#    zh-cn / de-de sheets.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
